# Re-order the comma-separated list of "Recorded By" values in column G
# so that they are sorted alphabetically in a case-insensitive manner
# (stable sort, so entries that compare equal keep their relative order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $sorted = $trimmed | Sort-Object

        $newVal = [string]::Join(", ", $sorted)
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
